# Auto-generated Excel COM-interop edit script
# Applies the scheduled-runner price/profit recalculation updates to the
# Leve profit worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 330.7143
$ws.Range("I33").Value = 219.57143
$ws.Range("J33").Value = 775.2857
$ws.Range("K33").Value = 219.57143
$ws.Range("L33").Value = 775.2857
$ws.Range("M33").Value = 9.428570000000008
$ws.Range("N33").Value = -1233.2857

$ws.Range("H103").Value = 1424.8
$ws.Range("I103").Value = 2000
$ws.Range("J103").Value = 1041.3334
$ws.Range("K103").Value = 6000
$ws.Range("L103").Value = 3124.0002
$ws.Range("M103").Value = -5414
$ws.Range("N103").Value = -4296.0002

$ws.Range("H116").Value = 56975
$ws.Range("I116").Value = 81462.5
$ws.Range("K116").Value = 81462.5
$ws.Range("M116").Value = -78020.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7442212
$ws.Range("I32").Value = 3402365
$ws.Range("K32").Value = 3402365
$ws.Range("M32").Value = -3402078

$ws.Range("H122").Value = 6471.288
$ws.Range("I122").Value = 5238.8535
$ws.Range("K122").Value = 15716.5605
$ws.Range("M122").Value = -13266.5605

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 50180.934
$ws.Range("J62").Value = 50180.934
$ws.Range("L62").Value = 50180.934
$ws.Range("N62").Value = -51552.934

$ws.Range("H65").Value = 50180.934
$ws.Range("J65").Value = 50180.934
$ws.Range("L65").Value = 150542.802
$ws.Range("N65").Value = -157406.802

$ws.Range("H117").Value = 139994.5
$ws.Range("J117").Value = 139994.5
$ws.Range("L117").Value = 139994.5
$ws.Range("N117").Value = -149172.5

$ws.Range("H123").Value = 65000
$ws.Range("J123").Value = 65000
$ws.Range("L123").Value = 65000
$ws.Range("N123").Value = -74800

$ws.Range("H129").Value = 157387
$ws.Range("J129").Value = 157387
$ws.Range("L129").Value = 157387
$ws.Range("N129").Value = -167387

$ws.Range("H134").Value = 21307740
$ws.Range("I134").Value = 4204337
$ws.Range("K134").Value = 12613011
$ws.Range("M134").Value = -12610476

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 312.5
$ws.Range("I22").Value = 312.5
$ws.Range("K22").Value = 312.5
$ws.Range("M22").Value = 37.5

$ws.Range("H31").Value = 2956.037
$ws.Range("I31").Value = 2125.7144
$ws.Range("K31").Value = 2125.7144
$ws.Range("M31").Value = -1830.7144

$ws.Range("H34").Value = 2956.037
$ws.Range("I34").Value = 2125.7144
$ws.Range("K34").Value = 2125.7144
$ws.Range("M34").Value = -1923.7144

$ws.Range("H105").Value = 2426.3333
$ws.Range("I105").Value = 1889.75
$ws.Range("J105").Value = 3499.5
$ws.Range("K105").Value = 1889.75
$ws.Range("L105").Value = 3499.5
$ws.Range("M105").Value = -142.75
$ws.Range("N105").Value = -6993.5

$ws.Range("H134").Value = 2296.2334
$ws.Range("I134").Value = 1690.65
$ws.Range("K134").Value = 5071.950000000001
$ws.Range("M134").Value = -2536.950000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 733.3333
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""

$ws.Range("H84").Value = 733.3333
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""

$ws.Range("H98").Value = 1027
$ws.Range("J98").Value = 1058
$ws.Range("L98").Value = 3174
$ws.Range("N98").Value = -6170

$ws.Range("H107").Value = 688.89655
$ws.Range("J107").Value = 657.5
$ws.Range("L107").Value = 1972.5
$ws.Range("N107").Value = -5812.5

$ws.Range("H129").Value = 1644.3334
$ws.Range("I129").Value = 800
$ws.Range("J129").Value = 2066.5
$ws.Range("K129").Value = 2400
$ws.Range("L129").Value = 6199.5
$ws.Range("M129").Value = 2600
$ws.Range("N129").Value = -16199.5

$ws.Range("H140").Value = 40001270
$ws.Range("I140").Value = 40001270
$ws.Range("K140").Value = 120003810
$ws.Range("M140").Value = -119998630

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12905.4
$ws.Range("J70").Value = 4499.905
$ws.Range("L70").Value = 4499.905
$ws.Range("N70").Value = -5039.905

$ws.Range("H73").Value = 12905.4
$ws.Range("J73").Value = 4499.905
$ws.Range("L73").Value = 4499.905
$ws.Range("N73").Value = -6371.905

$ws.Range("H80").Value = 1888.4117
$ws.Range("I80").Value = 1167.1666
$ws.Range("J80").Value = 2281.818
$ws.Range("K80").Value = 1167.1666
$ws.Range("L80").Value = 2281.818
$ws.Range("M80").Value = -169.1666
$ws.Range("N80").Value = -4277.818

$ws.Range("H83").Value = 1888.4117
$ws.Range("I83").Value = 1167.1666
$ws.Range("J83").Value = 2281.818
$ws.Range("K83").Value = 5835.833000000001
$ws.Range("L83").Value = 11409.09
$ws.Range("M83").Value = -843.8330000000005
$ws.Range("N83").Value = -21393.09

$ws.Range("H107").Value = 1115.52
$ws.Range("I107").Value = 1079.7059
$ws.Range("J107").Value = 1191.625
$ws.Range("K107").Value = 1079.7059
$ws.Range("L107").Value = 1191.625
$ws.Range("M107").Value = 840.2941000000001
$ws.Range("N107").Value = -5031.625

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5333.3335
$ws.Range("I22").Value = 5000
$ws.Range("K22").Value = 5000
$ws.Range("M22").Value = -4705

$ws.Range("H27").Value = 5333.3335
$ws.Range("I27").Value = 5000
$ws.Range("K27").Value = 5000
$ws.Range("M27").Value = -4893

$ws.Range("H46").Value = 2993.318
$ws.Range("J46").Value = 4111.533
$ws.Range("L46").Value = 4111.533
$ws.Range("N46").Value = -4487.533

$ws.Range("H133").Value = 49886.5
$ws.Range("J133").Value = 49886.5
$ws.Range("L133").Value = 49886.5
$ws.Range("N133").Value = -54946.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 9999
$ws.Range("J29").Value = 9999
$ws.Range("L29").Value = 9999
$ws.Range("N29").Value = -10579

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = ""
$ws.Range("N40").Value = ""

$ws.Range("H122").Value = 9336.25

$ws.Range("H132").Value = 3069.389
$ws.Range("I132").Value = 2585.0312
$ws.Range("K132").Value = 7755.0936
$ws.Range("M132").Value = -5225.0936

$ws.Range("H136").Value = 2742.4
$ws.Range("I136").Value = 2178
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 6534
$ws.Range("M136").Value = -3984
$ws.Range("N136").Value = -20100

